# Commit: "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The "Periodo Mora" / "Valor Mora" (F) data block (rows 16-22) was refreshed:
# the period list got reversed in the database, and the "Salario Basico"
# (G) column was unified to a single new value for all workers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Periodo Mora" values (column E) for rows 16-22 -- old EC periods removed,
# new ones added (database refresh), ending up reversed relative to before.
$periods = @("2108", "2107", "2106", "2102", "2101", "2010", "2008")

# New "Valor Mora" values (column F) for rows 16-22, following the same reorder.
$valorMora = @(33942, 36341, 36341, 35112, 35112, 35112, 35112)

# New "Salario Basico" (column G) value, now identical for every row.
$salarioBasico = 908526

for ($i = 0; $i -lt 7; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valorMora[$i]
    $ws.Range("G$row").Value = $salarioBasico
}

# Column widths were nudged slightly wider after the database refresh
# (best-fit recalculated by Excel on save).
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.333333333333334
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333332
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
